# Regenerate the save_data "K" column (column G) values.
# The commit replaces the previously-stored "Strike#" values in column G
# with freshly calculated K values (std/mean based s_vals), row by row,
# for rows 2 through 33 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = [ordered]@{
    2  = 2
    3  = 3
    4  = 5
    5  = 3
    6  = 7
    7  = 5
    8  = 5
    9  = 7
    10 = 7
    11 = 7
    12 = 2
    13 = 9
    14 = 6
    15 = 6
    16 = 9
    17 = 5
    18 = 4
    19 = 6
    20 = 6
    21 = 5
    22 = 4
    23 = 7
    24 = 5
    25 = 6
    26 = 5
    27 = 7
    28 = 6
    29 = 4
    30 = 4
    31 = 5
    32 = 5
    33 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
